$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.466.76"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.970.02"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'377.99"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").Value = "'104.70"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").Value = "'0.540"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").Value = "'37.17"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "'0.0842"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "3.435.53"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "'18.40"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "'7.54"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").Value = "2.970.01"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").Value = "'0.965"
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "51.398.83"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D20").Value = "'7.38"
$ws.Range("E20").Value = "  +2.85%  "
$ws.Range("D21").Value = "'12.89"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").Value = "'69.39"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "'261.38"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "'2.82"
$ws.Range("E25").Value = "  +4.73%  "
$ws.Range("B26").Value = "Filecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D26").Value = "'8.08"
$ws.Range("E26").Value = "  +16.64%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'7.58"
$ws.Range("E27").Value = "  +23.42%  "
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.112"
$ws.Range("E30").Value = "  +8.70%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'25.84"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "'35.00"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "'50.89"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "'0.0445"
$ws.Range("E36").Value = "  +5.77%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'3.03"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "'17.17"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").Value = "'2.59"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "'0.115"
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("D43").Value = "'124.70"
$ws.Range("E43").Value = "  +4.48%  "
$ws.Range("D44").Value = "'21.67"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("D45").Value = "'0.289"
$ws.Range("E45").Value = "  +19.00%  "
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("D47").Value = "'2.36"
$ws.Range("E47").Value = "  +2.70%  "
$ws.Range("D48").Value = "2.033.92"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").Value = "'3.21"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("D50").Value = "'0.0340"
$ws.Range("E50").Value = "  +10.23%  "
$ws.Range("E51").Value = "  +2.58%  "
